$d = $word.ActiveDocument

function Set-RunText {
    param(
        [int]$ParaIndex,
        [string[]]$NewTexts
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    # the bullet paragraphs always begin with a 1-character run
    # (the opening curly quote) that must be left untouched; only the
    # text after it (up to, but excluding, the paragraph mark) is replaced.
    $target = $d.Range($start + 1, $end - 1)

    $runsXml = ""
    foreach ($t in $NewTexts) {
        $escaped = $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += "<w:r><w:rPr></w:rPr><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    }

    $xml = "<pkg:package xmlns:pkg='http://schemas.openxmlformats.org/package/2006/metadata/core-properties'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:body><w:p>$runsXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

    $target.InsertXML($xml)
}

# 1) "EventID" : "nnnn" contains the Windows event id  ->  "_source_tag": "windows_agent"
Set-RunText 59 @("_source_tag”: “windows_agent”")

# 2) "EventLog": "xxx" ... -> "event_id": "nnnn" contains the Windows event id
Set-RunText 60 @("event_id”: “nnnn” contains the Windows event id")

# 3) "_source_type" : "WindowsAgent" ... -> "event_log":  "xxx" ... contains the name of the event log that produced the message
Set-RunText 61 @("event_log”:  ”xxx” … contains the name of the event log that produced the message")

# 4) "_log_type": "eventlog" OR "_log_type": "file" ... -> split across two runs: "l" + "og_type": ...
Set-RunText 62 @("l", "og_type”: ”eventlog” OR “_log_type”: ”file” … indicates whether the log message originated in a Windows event log or originated from the “tail” operation")

# 5) Give the (empty, text-less) run right before the page break in the
#    paragraph that follows "Log File Name" an explicit sz/szCs of 28,
#    matching the paragraph mark's run formatting, while leaving the
#    pPr and the page-break run untouched.
$pageBreakPara = $d.Paragraphs.Item(68)
$pbRange = $pageBreakPara.Range
$xml = "<pkg:package xmlns:pkg='http://schemas.openxmlformats.org/package/2006/metadata/core-properties'>" +
       "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
       "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:body><w:p><w:pPr><w:pStyle w:val=`"Normal`"/><w:rPr><w:rStyle w:val=`"Strong`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr>" +
       "<w:r><w:rPr><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:r><w:r><w:br w:type=`"page`"/></w:r></w:p>" +
       "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$pbRange.InsertXML($xml)

# 6) The paragraph containing "Other" gains an explicit spacing
#    override (before=0, after=160 twips = 8pt), matching what was
#    already the effective inherited spacing.
$otherPara = $d.Paragraphs.Item(69)
$otherPara.Range.ParagraphFormat.SpaceBefore = 0
$otherPara.Range.ParagraphFormat.SpaceAfter = 8

Write-Host "done"
